$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

# Update existing values
$ws.Range("B2").Value = "ctc"
$ws.Range("B3").Value = "levy"
$ws.Range("A4").Value = "Hello it is two twenty"

# Add new values
$ws.Range("B4").Value = "child"
$ws.Range("A5").Value = "A5"
$ws.Range("B5").Value = "credit"

# Update selection to D7
$ws.Range("D7").Select()
